$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 36822.953
$ws.Range("I98").Value = 1098.5333
$ws.Range("J98").Value = 113375.29
$ws.Range("K98").Value = 1098.5333
$ws.Range("L98").Value = 113375.29
$ws.Range("M98").Value = 399.4666999999999
$ws.Range("N98").Value = -116371.29
$ws.Range("H112").Value = 1757.3334
$ws.Range("I112").Value = 150
$ws.Range("J112").Value = 1827.2174
$ws.Range("K112").Value = 450
$ws.Range("L112").Value = 5481.6522
$ws.Range("M112").Value = 658
$ws.Range("N112").Value = -7697.6522
$ws.Range("H113").Value = 2875
$ws.Range("I113").Value = 2800
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 2800
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = 454
$ws.Range("N113").Value = -9408
$ws.Range("H122").Value = 36822.953
$ws.Range("I122").Value = 1098.5333
$ws.Range("J122").Value = 113375.29
$ws.Range("K122").Value = 3295.5999
$ws.Range("L122").Value = 340125.87
$ws.Range("M122").Value = -845.5999000000002
$ws.Range("N122").Value = -345025.87
$ws.Range("H125").Value = 1354.7273
$ws.Range("I125").Value = 1251
$ws.Range("K125").Value = 11259
$ws.Range("M125").Value = -8799
$ws.Range("H132").Value = 30671.818
$ws.Range("I132").Value = 18993.678
$ws.Range("J132").Value = 129101.86
$ws.Range("K132").Value = 56981.034
$ws.Range("L132").Value = 387305.58
$ws.Range("M132").Value = -54451.034
$ws.Range("N132").Value = -392365.58

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1707.1154
$ws.Range("I2").Value = 1615.4
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 1615.4
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -1502.4
$ws.Range("N2").Value = -4226
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H43").Value = 1006924.7
$ws.Range("J43").Value = 1006924.7
$ws.Range("L43").Value = 1006924.7
$ws.Range("N43").Value = -1007550.7
$ws.Range("H45").Value = 35716104
$ws.Range("I45").Value = 47620610
$ws.Range("J45").Value = 2589.2856
$ws.Range("K45").Value = 47620610
$ws.Range("L45").Value = 2589.2856
$ws.Range("M45").Value = -47620233
$ws.Range("N45").Value = -3343.2856
$ws.Range("H74").Value = 1452.4736
$ws.Range("I74").Value = 1187.5294
$ws.Range("K74").Value = 1187.5294
$ws.Range("M74").Value = -313.5293999999999
$ws.Range("H77").Value = 1452.4736
$ws.Range("I77").Value = 1187.5294
$ws.Range("K77").Value = 5937.646999999999
$ws.Range("M77").Value = -1569.646999999999
$ws.Range("H116").Value = 1707.1154
$ws.Range("I116").Value = 1615.4
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 1615.4
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 678.5999999999999
$ws.Range("N116").Value = -8588
$ws.Range("H122").Value = 2302.25
$ws.Range("I122").Value = 1976
$ws.Range("J122").Value = 2628.5
$ws.Range("K122").Value = 5928
$ws.Range("L122").Value = 7885.5
$ws.Range("M122").Value = -3478
$ws.Range("N122").Value = -12785.5
$ws.Range("H135").Value = 48054.832
$ws.Range("J135").Value = 48054.832
$ws.Range("L135").Value = 48054.832
$ws.Range("N135").Value = -58194.832

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1707.1154
$ws.Range("I3").Value = 1615.4
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 1615.4
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -1501.4
$ws.Range("N3").Value = -4228
$ws.Range("H105").Value = 2005.4286
$ws.Range("I105").Value = 1829.5883
$ws.Range("J105").Value = 2752.75
$ws.Range("K105").Value = 1829.5883
$ws.Range("L105").Value = 2752.75
$ws.Range("M105").Value = -82.58829999999989
$ws.Range("N105").Value = -6246.75
$ws.Range("H134").Value = 2209.6667
$ws.Range("I134").Value = 1814.5588
$ws.Range("J134").Value = 3888.875
$ws.Range("K134").Value = 5443.6764
$ws.Range("L134").Value = 11666.625
$ws.Range("M134").Value = -2908.6764
$ws.Range("N134").Value = -16736.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 26424382
$ws.Range("I4").Value = 100001300
$ws.Range("J4").Value = 146911.78
$ws.Range("K4").Value = 100001300
$ws.Range("L4").Value = 146911.78
$ws.Range("M4").Value = -100001188
$ws.Range("N4").Value = -147135.78
$ws.Range("H58").Value = 1737.0883
$ws.Range("I58").Value = 1166.2727
$ws.Range("J58").Value = 2783.5833
$ws.Range("K58").Value = 1166.2727
$ws.Range("L58").Value = 2783.5833
$ws.Range("M58").Value = -963.2727
$ws.Range("N58").Value = -3189.5833
$ws.Range("H132").Value = 426914.66
$ws.Range("I132").Value = 1733.6957
$ws.Range("J132").Value = 1404830.9
$ws.Range("K132").Value = 5201.0871
$ws.Range("L132").Value = 4214492.699999999
$ws.Range("M132").Value = -2671.0871
$ws.Range("N132").Value = -4219552.699999999
$ws.Range("H136").Value = 1737.0883
$ws.Range("I136").Value = 1166.2727
$ws.Range("J136").Value = 2783.5833
$ws.Range("K136").Value = 3498.8181
$ws.Range("L136").Value = 8350.749899999999
$ws.Range("M136").Value = -948.8181
$ws.Range("N136").Value = -13450.7499

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 747.5
$ws.Range("I4").Value = 96.666664
$ws.Range("J4").Value = 2700
$ws.Range("K4").Value = 289.999992
$ws.Range("L4").Value = 8100
$ws.Range("M4").Value = -177.999992
$ws.Range("N4").Value = -8324
$ws.Range("H39").Value = 2266.6667
$ws.Range("J39").Value = 2454.5454
$ws.Range("L39").Value = 7363.6362
$ws.Range("N39").Value = -7951.6362
$ws.Range("H68").Value = 5920986
$ws.Range("I68").Value = 4445174
$ws.Range("J68").Value = 6945855
$ws.Range("K68").Value = 13335522
$ws.Range("L68").Value = 20837565
$ws.Range("M68").Value = -13334711
$ws.Range("N68").Value = -20839187
$ws.Range("H71").Value = 5920986
$ws.Range("I71").Value = 4445174
$ws.Range("J71").Value = 6945855
$ws.Range("K71").Value = 40006566
$ws.Range("L71").Value = 62512695
$ws.Range("M71").Value = -40002510
$ws.Range("N71").Value = -62520807
$ws.Range("H80").Value = 31330944
$ws.Range("J80").Value = 32277104
$ws.Range("L80").Value = 96831312
$ws.Range("N80").Value = -96833184
$ws.Range("H83").Value = 31330944
$ws.Range("J83").Value = 32277104
$ws.Range("L83").Value = 290493936
$ws.Range("N83").Value = -290503296
$ws.Range("H92").Value = 8447932
$ws.Range("I92").Value = 19609176
$ws.Range("J92").Value = 76998.25
$ws.Range("K92").Value = 58827528
$ws.Range("L92").Value = 230994.75
$ws.Range("M92").Value = -58826280
$ws.Range("N92").Value = -233490.75
$ws.Range("H98").Value = 340.2
$ws.Range("I98").Value = 203
$ws.Range("J98").Value = 374.5
$ws.Range("K98").Value = 609
$ws.Range("L98").Value = 1123.5
$ws.Range("M98").Value = 889
$ws.Range("N98").Value = -4119.5
$ws.Range("H129").Value = 121528.32
$ws.Range("I129").Value = 375861.25
$ws.Range("J129").Value = 1842.2354
$ws.Range("K129").Value = 1127583.75
$ws.Range("L129").Value = 5526.706200000001
$ws.Range("M129").Value = -1122583.75
$ws.Range("N129").Value = -15526.7062

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 30
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 30
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 82
$ws.Range("N5").ClearContents()
$ws.Range("H116").Value = 42742
$ws.Range("J116").Value = 42742
$ws.Range("L116").Value = 42742
$ws.Range("N116").Value = -51920
$ws.Range("H122").Value = 1456.2858
$ws.Range("I122").Value = 1372.4286
$ws.Range("J122").Value = 1540.1428
$ws.Range("K122").Value = 4117.2858
$ws.Range("L122").Value = 4620.428400000001
$ws.Range("M122").Value = -1667.2858
$ws.Range("N122").Value = -9520.428400000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 7208.3335
$ws.Range("H7").Value = 45457380
$ws.Range("I7").Value = 100001850
$ws.Range("J7").Value = 3658.75
$ws.Range("K7").Value = 100001850
$ws.Range("L7").Value = 3658.75
$ws.Range("M7").Value = -100001738
$ws.Range("N7").Value = -3882.75
$ws.Range("H16").Value = 3150.524
$ws.Range("I16").Value = 3064.7778
$ws.Range("J16").Value = 3665
$ws.Range("K16").Value = 3064.7778
$ws.Range("L16").Value = 3665
$ws.Range("M16").Value = -2894.7778
$ws.Range("N16").Value = -4005
$ws.Range("H17").Value = 2952.5
$ws.Range("I17").Value = 2952.5
$ws.Range("K17").Value = 2952.5
$ws.Range("M17").Value = -2782.5
$ws.Range("H61").Value = 2420
$ws.Range("I61").Value = 1400
$ws.Range("J61").Value = 2675
$ws.Range("K61").Value = 1400
$ws.Range("L61").Value = 2675
$ws.Range("M61").Value = -1198
$ws.Range("N61").Value = -3079
$ws.Range("H113").Value = 2420
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 2675
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 2675
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -7015
$ws.Range("H122").Value = 145257.72
$ws.Range("I122").Value = 334934.66
$ws.Range("K122").Value = 1004803.98
$ws.Range("M122").Value = -1002353.98
$ws.Range("H126").Value = 45457380
$ws.Range("I126").Value = 100001850
$ws.Range("J126").Value = 3658.75
$ws.Range("K126").Value = 300005550
$ws.Range("L126").Value = 10976.25
$ws.Range("M126").Value = -300003080
$ws.Range("N126").Value = -15916.25
$ws.Range("H136").Value = 1762.6216
$ws.Range("I136").Value = 1346.3549
$ws.Range("K136").Value = 4039.0647
$ws.Range("M136").Value = -1489.0647

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 10002
$ws.Range("I2").Value = 10002
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 10002
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -9890
$ws.Range("N2").ClearContents()
$ws.Range("H122").Value = 28571428
$ws.Range("I122").Value = 28571428
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 85714284
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -85711834
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1547.3125
$ws.Range("I132").Value = 1192.1086
$ws.Range("K132").Value = 3576.3258
$ws.Range("M132").Value = -1192.1086
